$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Valor Mora" total (E11) ---
$ws.Range("E11").Value2 = 172718

# --- 2. Update "Cant. Periodos" (F13) ---
$ws.Range("F13").Value2 = 4

# --- 3. Insert a new row for the new period (2508), pushing the signature block down ---
$ws.Rows.Item(19).Insert()

# Give the new row the same bordered-table look as the rest of the data rows
$rng19 = $ws.Range("B19:J19")
$rng19.Borders.Item(7).LineStyle = 1    # xlEdgeLeft
$rng19.Borders.Item(7).Weight = 2
$rng19.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$rng19.Borders.Item(8).Weight = 2
$rng19.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$rng19.Borders.Item(9).Weight = 2
$rng19.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$rng19.Borders.Item(10).Weight = 2
$rng19.Borders.Item(11).LineStyle = 1   # xlInsideVertical
$rng19.Borders.Item(11).Weight = 2

# --- 4. Refresh the period rows (16-19) into chronological order: 2505, 2506, 2507, 2508 ---
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "1002293755"
$ws.Range("D16").Value2 = "LUIS JOSE BARRIOS GUZMAN"
$ws.Range("E16").Value2 = "2505"
$ws.Range("F16").Value2 = 1898
$ws.Range("G16").Value2 = 1423500

$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "1002293755"
$ws.Range("D17").Value2 = "LUIS JOSE BARRIOS GUZMAN"
$ws.Range("E17").Value2 = "2506"
$ws.Range("F17").Value2 = 56940
$ws.Range("G17").Value2 = 1423500

$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1002293755"
$ws.Range("D18").Value2 = "LUIS JOSE BARRIOS GUZMAN"
$ws.Range("E18").Value2 = "2507"
$ws.Range("F18").Value2 = 56940
$ws.Range("G18").Value2 = 1423500

$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1002293755"
$ws.Range("D19").Value2 = "LUIS JOSE BARRIOS GUZMAN"
$ws.Range("E19").Value2 = "2508"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500
